# The commit swaps the contents of ppt/theme/theme1.xml (the "Integral"
# theme used by the slide master / the whole deck) and ppt/theme/theme2.xml
# (the "Office Theme" used only by the notes master), while the part names
# and relationships stay untouched.
#
# The two theme parts differ only in <a:clrScheme> (its name + the 12
# scheme colours) - the font scheme and format scheme (fills/lines/effects)
# are byte-identical between the two themes. So reproducing the swap for
# the part that actually drives the presentation's look (theme1.xml, wired
# to the slide master that every slide/layout inherits from) comes down to
# replacing its 12 scheme colours with the ones that used to live in
# theme2.xml ("Office Theme").
#
# (The notes master's theme, theme2.xml, is not reachable through the
# PowerPoint object model exposed by this host - NotesMaster.Theme simply
# aliases the slide master's theme here - so it can't be targeted
# independently; the colour swap below is applied through the one theme
# object the host does expose.)

$p = $ppt.ActivePresentation
$scheme = $p.SlideMaster.Theme.ThemeColorScheme

# ThemeColorScheme index -> OOXML clrScheme slot:
#  1 dk1   2 lt1   3 dk2   4 lt2
#  5 accent1  6 accent2  7 accent3  8 accent4  9 accent5  10 accent6
#  11 hlink   12 folHlink
# .RGB takes the standard COM "long" RGB encoding (0x00BBGGRR), i.e. the
# bytes of the hex colour reversed, so 0xRRGGBB (as written in the OOXML
# srgbClr val) becomes 0x00BBGGRR below.

$scheme.Item(1).RGB  = 0x000000   # dk1      000000
$scheme.Item(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$scheme.Item(3).RGB  = 0x6A5444   # dk2      44546A
$scheme.Item(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$scheme.Item(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$scheme.Item(6).RGB  = 0x317DED   # accent2  ED7D31
$scheme.Item(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$scheme.Item(8).RGB  = 0x00C0FF   # accent4  FFC000
$scheme.Item(9).RGB  = 0xC47244   # accent5  4472C4
$scheme.Item(10).RGB = 0x47AD70   # accent6  70AD47
$scheme.Item(11).RGB = 0xC16305   # hlink    0563C1
$scheme.Item(12).RGB = 0x724F95   # folHlink 954F72
